## Auto commit at 2025-09-15 7:55:36.53
## Appends the new daily readings for 2025-09-14 (Excel serial 45914) for
## both charging stations, extending the log by two rows (90 and 91),
## then updates the sheet's scroll position / selection to where the user
## ended up after entering the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 90: 四方坪站 (Sifangping station), 2025-09-14 ----
$ws.Range("A90").Value = 45914
$ws.Range("B90").Value = "四方坪站"
$ws.Range("C90").Value = 10563.62
$ws.Range("D90").Value = 8536.7199999999993
$ws.Range("E90").Value = 3622.01
$ws.Range("F90").Value = 429

# ---- Row 91: 高岭站 (Gaoling station), 2025-09-14 ----
$ws.Range("A91").Value = 45914
$ws.Range("B91").Value = "高岭站"
$ws.Range("C91").Value = 4957.55
$ws.Range("D91").Value = 3716.93
$ws.Range("E91").Value = 1255.78
$ws.Range("F91").Value = 165

# The used range now extends to F91 automatically; scroll the window down
# so row 82 is at the top and put the active selection on G93, matching
# where the cursor ended up after the new rows were typed in.
$win = $excel.ActiveWindow
$win.ScrollRow = 82
$win.ScrollColumn = 1
$ws.Range("G93").Select()
